$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023xxxx_listplayer")

$ws.Range("F4").Value = 200
$ws.Range("F6").Value = 81
$ws.Range("F7").Value = 45
$ws.Range("F8").Value = 200
$ws.Range("F12").Value = 36
$ws.Range("F37").Value = 36
$ws.Range("F57").Value = 200
$ws.Range("F64").Value = 35
$ws.Range("F67").Value = 68
